$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = -0.804
$ws.Range("E2").Value = -0.146
$ws.Range("I2").Value = -48.43137254901961
$ws.Range("J2").Value = -48.43137254901961
$ws.Range("K2").Value = 68.5
$ws.Range("L2").Value = 1343.137254901961
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 69.8
$ws.Range("V2").Value = 0.02083333333333333
$ws.Range("W2").Value = 0.02193052665279334
$ws.Range("X2").Value = 0.01872458975801795
$ws.Range("Y2").Value = 0.003205936894775389
$ws.Range("Z2").Value = 0.00001451977964099133
$ws.Range("AA2").Value = -0.0007032128571225213
$ws.Range("AB2").Value = 0.01859028329849068
$ws.Range("AC2").Value = -0.0192934961556132
$ws.Range("AD2").Value = 386
$ws.Range("AF2").Value = 386
$ws.Range("AG2").Value = 316.2
$ws.Range("AH2").Value = 0.1033079970024623
$ws.Range("AI2").Value = 0.1089318470438832
$ws.Range("AJ2").Value = 0.08623793159875634
$ws.Range("AK2").Value = 0.09102685896882287
$ws.Range("AL2").Value = 16.3
$ws.Range("AM2").Value = -70.9
$ws.Range("AO2").Value = -0.1515337423312884
$ws.Range("AQ2").Value = 0.03483779971791255

# Row 2: clear buybacks_cash_returned (T2) cell entirely
$ws.Range("T2").ClearContents()

# Row 3 updates
$ws.Range("D3").Value = -0.804
$ws.Range("E3").Value = -0.146
$ws.Range("I3").Value = -48.43137254901961
$ws.Range("J3").Value = -48.43137254901961
$ws.Range("K3").Value = 68.5
$ws.Range("L3").Value = 1343.137254901961
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 69.8
$ws.Range("V3").Value = 0.02083333333333333
$ws.Range("W3").Value = 0.02193052665279334
$ws.Range("X3").Value = 0.01872458975801795
$ws.Range("Y3").Value = 0.003205936894775389
$ws.Range("Z3").Value = 0.00001451977964099133
$ws.Range("AA3").Value = -0.0007032128571225213
$ws.Range("AB3").Value = 0.01859028329849068
$ws.Range("AC3").Value = -0.0192934961556132
$ws.Range("AD3").Value = 386
$ws.Range("AF3").Value = 386
$ws.Range("AG3").Value = 316.2
$ws.Range("AH3").Value = 0.1033079970024623
$ws.Range("AI3").Value = 0.1089318470438832
$ws.Range("AJ3").Value = 0.08623793159875634
$ws.Range("AK3").Value = 0.09102685896882287
$ws.Range("AL3").Value = 16.3
$ws.Range("AM3").Value = -70.9
$ws.Range("AO3").Value = -0.1515337423312884
$ws.Range("AQ3").Value = 0.03483779971791255

# Row 3: clear buybacks_cash_returned (T3) cell entirely
$ws.Range("T3").ClearContents()

Write-Host "Capital structure database updated"
